# Apply the edit described by the diff:
# On the "Usuarios" worksheet, fill column B ("Incluido") with "X" for the
# remaining rows (15-18) that were previously left blank, and leave the
# final selection on B18 (matching the new <selection activeCell="B18" sqref="B18"/>).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Usuarios")

$ws.Range("B15").Value = "X"
$ws.Range("B16").Value = "X"
$ws.Range("B17").Value = "X"
$ws.Range("B18").Value = "X"

$ws.Activate()
$ws.Range("B18").Select()
